$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "67.492.33"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "  -1.00%  "

$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "3.759.83"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "  -2.25%  "

$ws.Range("E4").Value = "  +0.31%  "

$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "596.16"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "  -0.81%  "

$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "168.59"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "  +0.34%  "

$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "3.757.63"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "  -2.33%  "

$ws.Range("E8").Value = "  +0.09%  "

$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.529"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "  +0.01%  "

$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.165"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "  +0.93%  "

$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "6.49"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "  +0.65%  "

$ws.Range("E12").Value = "  -0.09%  "

$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.0000272"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "  +4.86%  "

$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "36.78"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "  -0.48%  "

$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "4.390.43"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "  -2.04%  "

$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "3.766.52"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "  -1.73%  "

$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "18.98"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "  +5.25%  "

$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "67.542.60"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "  -0.89%  "

$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "7.25"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "  -1.10%  "

$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "10.53"
$ws.Range("D21").Style = "Normal"
$ws.Range("E21").Value = "  -1.87%  "

$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "467.21"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "  +0.08%  "

$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = "0.723"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "  -1.33%  "

$ws.Range("E24").Value = "  -6.15%  "

$ws.Range("E25").Value = "  +1.19%  "

$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "2.24"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "  +0.52%  "

$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "12.14"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "  +0.83%  "

$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = "10.30"
$ws.Range("D28").Style = "Normal"
$ws.Range("E28").Value = "  +3.05%  "

$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "0.999"
$ws.Range("D29").Style = "Normal"
$ws.Range("E29").Value = "  +0.05%  "

$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "2.90"
$ws.Range("D30").Style = "Normal"
$ws.Range("E30").Value = "  -2.13%  "

$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "3.919.79"
$ws.Range("D31").Style = "Normal"
$ws.Range("E31").Value = "  -1.70%  "

$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "7.63"
$ws.Range("D32").Style = "Normal"
$ws.Range("E32").Value = "  +0.46%  "

$ws.Range("E33").Value = "  -2.51%  "

$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "30.40"
$ws.Range("D34").Style = "Normal"
$ws.Range("E34").Value = "  -2.59%  "

$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "9.15"
$ws.Range("D35").Style = "Normal"
$ws.Range("E35").Value = "  -3.47%  "

$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "3.721.91"
$ws.Range("D36").Style = "Normal"
$ws.Range("E36").Value = "  -2.16%  "

$ws.Range("E37").Value = "  +5.24%  "

$ws.Range("E38").Value = "  +0.80%  "

$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = "5.89"
$ws.Range("D39").Style = "Normal"
$ws.Range("E39").Value = "  -0.15%  "

$ws.Range("E40").Value = "  -2.22%  "

$ws.Range("E41").Value = "  -2.38%  "

$ws.Range("E42").Value = "  +0.43%  "

$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.315"
$ws.Range("D43").Style = "Normal"
$ws.Range("E43").Value = "  +0.79%  "

$ws.Range("E44").Value = "  +0.01%  "

$ws.Range("B45").Value = "Cosmos"
$ws.Range("C45").Value = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.71"
$ws.Range("D45").Style = "Normal"
$ws.Range("E45").Value = "  +1.55%  "

$ws.Range("B46").Value = "Stacks"
$ws.Range("C46").Value = "https://coinranking.com/coin/mMPrMcB7+stacks-stx"
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "1.96"
$ws.Range("D46").Style = "Normal"
$ws.Range("E46").Value = "  -0.44%  "

$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "46.26"
$ws.Range("D47").Style = "Normal"
$ws.Range("E47").Value = "  -1.68%  "

$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "401.93"
$ws.Range("D48").Style = "Normal"
$ws.Range("E48").Value = "  -4.24%  "

$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "0.000276"
$ws.Range("D49").Style = "Normal"
$ws.Range("E49").Value = "  -6.58%  "

$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "141.56"
$ws.Range("D50").Style = "Normal"
$ws.Range("E50").Value = "  -1.04%  "

$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.0355"
$ws.Range("D51").Style = "Normal"
$ws.Range("E51").Value = "  -0.39%  "
